# Documentation.docx - "Added comments for files"
#
# Updates the second paragraph of the document (the one starting "For
# creating the front-end for the website ...") so that it:
#
#   1. Mentions that a navbar was created with materialize for each page,
#      right after "...website's design easy to customize" and before
#      the following ". The functionality ..." sentence.
#   2. Drops the "bcryptjs" mention, keeping only "jsonwebtoken", and
#      moves the "where the jwt-token" clause right after it:
#        "...jsonwebtoken- and bcryptjs-libraries, where the jwt-token..."
#        -> "...jsonwebtoken, where the jwt-token..."
#   3. Rewrites the closing sentence: the jwt-token is now said to be
#      kept in the browser's cookies (instead of local storage), and
#      those cookies are checked when the token is needed for
#      authentication (instead of "local storage ... each time
#      authentication is needed"):
#        "...The jwt-token is stored in the browser's local storage in
#        cookies and the local storage is checked for the token each
#        time authentication is needed."
#        -> "...The jwt-token is stored in the browser's cookies and the
#        cookies are checked when the token is needed for
#        authentication."

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll   = 2
$wdCollapseEnd  = 0

# --- 1) Insert the new "navbar" clause -------------------------------
$r = $d.Content
$r.Find.Execute("design easy to customize", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "", 0)
$r.Collapse($wdCollapseEnd)
$r.InsertAfter(" and I created a navbar using materialize for each page")

# --- 2) Remove the bcryptjs mention / reorder the jwt-token clause ---
$d.Content.Find.Execute("jsonwebtoken- and bcryptjs-libraries, where the jwt-token", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "jsonwebtoken, where the jwt-token", $wdReplaceAll)

# --- 3) Rewrite where/how the jwt-token is kept and checked ----------
$d.Content.Find.Execute("The jwt-token is stored in the browser's local storage in cookies and the local storage is checked for the token each time authentication is needed.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "The jwt-token is stored in the browser's cookies and the cookies are checked when the token is needed for authentication.", $wdReplaceAll)
